$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("separation")

$ws.Range("F3").Value = 0.05331039138987579
$ws.Range("I3").Value = 0.05198846281844722
$ws.Range("F4").Value = 0.03054070656352176
$ws.Range("G4").Value = 0.1453078904717736
$ws.Range("I4").Value = 0.029351635134950333
$ws.Range("J4").Value = 0.14411881904320215
$ws.Range("F5").Value = 0.025129467108006374
$ws.Range("G5").Value = 0.04659224249893037
$ws.Range("I5").Value = 0.024425324250863516
$ws.Range("J5").Value = 0.04553602821321608
$ws.Range("F6").Value = 0.07657326248711238
$ws.Range("G6").Value = 0.07040603735437492
$ws.Range("H6").Value = 0.18770442345298777
$ws.Range("I6").Value = 0.020924713544851106
$ws.Range("J6").Value = 0.023766180211517773
$ws.Range("K6").Value = 0.14106456631013062
$ws.Range("F7").Value = 0.108168586394391
$ws.Range("G7").Value = 0.06751898285526695
$ws.Range("H7").Value = 0.08954536183624093
$ws.Range("I7").Value = 0.017521549623078135
$ws.Range("J7").Value = 0.01879868831527512
$ws.Range("K7").Value = 0.0408250672962491
$ws.Range("F8").Value = 0.15541945590358547
$ws.Range("G8").Value = 0.1071633510318814
$ws.Range("H8").Value = 0.09942663073889114
$ws.Range("I8").Value = 0.017030007733732686
$ws.Range("J8").Value = 0.01783234479874739
$ws.Range("K8").Value = 0.020597554400399353
$ws.Range("F9").Value = 0.2287182902740183
$ws.Range("G9").Value = 0.18196677340452938
$ws.Range("H9").Value = 0.12805640472117413
$ws.Range("I9").Value = 0.01808979680878255
$ws.Range("J9").Value = 0.018286128977983478
$ws.Range("K9").Value = 0.01956326767018046
$ws.Range("F10").Value = 0.29040071964001773
$ws.Range("G10").Value = 0.2818455149101809
$ws.Range("H10").Value = 0.2146250053318345
$ws.Range("I10").Value = 0.01822849707683276
$ws.Range("J10").Value = 0.01785782211670517
$ws.Range("K10").Value = 0.01866015918171987
$ws.Range("F11").Value = 0.4194297881336389
$ws.Range("G11").Value = 0.4217402668669335
$ws.Range("H11").Value = 0.40178256944167096
$ws.Range("I11").Value = 0.015337469227748111
$ws.Range("J11").Value = 0.014727878299518448
$ws.Range("K11").Value = 0.014924210468719374
$ws.Range("F12").Value = 0.3898490913120772
$ws.Range("G12").Value = 0.39239436500549935
$ws.Range("H12").Value = 0.39506342309220926
$ws.Range("I12").Value = 0.015298202561081446
$ws.Range("J12").Value = 0.014923406592979372
$ws.Range("K12").Value = 0.01455273163285178
$ws.Range("F14").Value = 0.06408813614253714
$ws.Range("I14").Value = 0.06038673614253714
$ws.Range("F15").Value = 0.04790435614253713
$ws.Range("G15").Value = 0.18636255083386777
$ws.Range("I15").Value = 0.044574956142537134
$ws.Range("J15").Value = 0.18303315083386776
$ws.Range("F16").Value = 0.043417616142537134
$ws.Range("G16").Value = 0.05850321614253713
$ws.Range("I16").Value = 0.04144601614253714
$ws.Range("J16").Value = 0.05554581614253713
$ws.Range("F17").Value = 0.0369988228092038
$ws.Range("G17").Value = 0.03872655614253714
$ws.Range("H17").Value = 0.18050775083386777
$ws.Range("I17").Value = 0.03613495614253714
$ws.Range("J17").Value = 0.03613495614253714
$ws.Range("K17").Value = 0.17791615083386778
$ws.Range("F18").Value = 0.034664056142537134
$ws.Range("G18").Value = 0.034664056142537134
$ws.Range("H18").Value = 0.049694556142537136
$ws.Range("I18").Value = 0.028646056142537134
$ws.Range("J18").Value = 0.028646056142537134
$ws.Range("K18").Value = 0.04367655614253713
$ws.Range("F19").Value = 0.09256705614253714
$ws.Range("G19").Value = 0.09134038947587048
$ws.Range("H19").Value = 0.08888705614253714
$ws.Range("I19").Value = 0.027519056142537135
$ws.Range("J19").Value = 0.027519056142537135
$ws.Range("K19").Value = 0.027519056142537135
$ws.Range("F20").Value = 0.18717600434944429
$ws.Range("G20").Value = 0.1390936644604367
$ws.Range("H20").Value = 0.12692820290765325
$ws.Range("I20").Value = 0.028808956142537135
$ws.Range("J20").Value = 0.028808956142537135
$ws.Range("K20").Value = 0.028808956142537135
$ws.Range("F21").Value = 0.26581846360420314
$ws.Range("G21").Value = 0.24604007812135775
$ws.Range("H21").Value = 0.2270824095717721
$ws.Range("I21").Value = 0.029491156142537136
$ws.Range("J21").Value = 0.029491156142537136
$ws.Range("K21").Value = 0.029491156142537136
$ws.Range("F22").Value = 0.25346067561107133
$ws.Range("G22").Value = 0.26163687066333924
$ws.Range("H22").Value = 0.2577701212913055
$ws.Range("I22").Value = 0.027767156142537136
$ws.Range("J22").Value = 0.027767156142537136
$ws.Range("K22").Value = 0.027767156142537136
$ws.Range("F23").Value = 0.24631379050306437
$ws.Range("G23").Value = 0.2544899855553323
$ws.Range("H23").Value = 0.26116081255456247
$ws.Range("I23").Value = 0.027649356142537135
$ws.Range("J23").Value = 0.027649356142537135
$ws.Range("K23").Value = 0.027649356142537135
$ws.Range("F36").Value = 0.08156797005810723
$ws.Range("I36").Value = 0.08156797005810723
$ws.Range("F37").Value = 0.03447518934582117
$ws.Range("G37").Value = 0.1896281062957387
$ws.Range("I37").Value = 0.03447518934582117
$ws.Range("J37").Value = 0.1896281062957387
$ws.Range("F38").Value = 0.025672378321953206
$ws.Range("G38").Value = 0.0660806902090109
$ws.Range("I38").Value = 0.025672378321953206
$ws.Range("J38").Value = 0.0660806902090109
$ws.Range("F39").Value = 0.02165367300976149
$ws.Range("G39").Value = 0.028538073009761487
$ws.Range("H39").Value = 0.17227354809652407
$ws.Range("I39").Value = 0.02165367300976149
$ws.Range("J39").Value = 0.028538073009761487
$ws.Range("K39").Value = 0.17227354809652407
$ws.Range("F40").Value = 0.021997966421895603
$ws.Range("G40").Value = 0.025829382498486556
$ws.Range("H40").Value = 0.0647252194414085
$ws.Range("I40").Value = 0.019497966421895604
$ws.Range("J40").Value = 0.023329382498486553
$ws.Range("K40").Value = 0.06222521944140849
$ws.Range("F41").Value = 0.16262844096262752
$ws.Range("G41").Value = 0.03166045215767163
$ws.Range("H41").Value = 0.038316080962627515
$ws.Range("I41").Value = 0.018878440962627516
$ws.Range("J41").Value = 0.021285452157671633
$ws.Range("K41").Value = 0.027941080962627516
$ws.Range("F42").Value = 0.24647490865118427
$ws.Range("G42").Value = 0.057042061984771894
$ws.Range("H42").Value = 0.04994262699753305
$ws.Range("I42").Value = 0.021160434283810516
$ws.Range("J42").Value = 0.021749430791413293
$ws.Range("K42").Value = 0.025580846868004242
$ws.Range("F43").Value = 0.31896245298784665
$ws.Range("G43").Value = 0.3139275557670384
$ws.Range("H43").Value = 0.16364790284248304
$ws.Range("I43").Value = 0.020894335087961142
$ws.Range("J43").Value = 0.01978231020757837
$ws.Range("K43").Value = 0.022189321402622482
$ws.Range("F44").Value = 0.28520766696589384
$ws.Range("G44").Value = 0.2833788941812049
$ws.Range("H44").Value = 0.2336299274233674
$ws.Range("I44").Value = 0.0139452515407072
$ws.Range("J44").Value = 0.012116478756018206
$ws.Range("K44").Value = 0.012705475263620987
$ws.Range("F45").Value = 0.2665065178087697
$ws.Range("G45").Value = 0.2653821299044635
$ws.Range("H45").Value = 0.26657063693897437
$ws.Range("I45").Value = 0.0139452515407072
$ws.Range("J45").Value = 0.01282086363640098
$ws.Range("K45").Value = 0.011708838756018206
